$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force specific cells that would otherwise be auto-recognized as numbers
# to remain plain text, matching the source data (text-formatted price strings).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values
$ws.Range("D2").Value = '61.904.17'
$ws.Range("E2").Value = '  -2.26%  '
$ws.Range("D3").Value = '2.577.38'
$ws.Range("E3").Value = '  -3.93%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '549.68'
$ws.Range("E5").Value = '  -0.68%  '
$ws.Range("D6").Value = '154.62'
$ws.Range("E6").Value = '  -2.12%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  +1.81%  '
$ws.Range("E9").Value = '  -1.36%  '
$ws.Range("E10").Value = '  -1.47%  '
$ws.Range("D11").Value = '5.55'
$ws.Range("E11").Value = '  +3.95%  '
$ws.Range("D12").Value = '0.364'
$ws.Range("E12").Value = '  -0.90%  '
$ws.Range("D13").Value = '3.036.78'
$ws.Range("E13").Value = '  -3.85%  '
$ws.Range("D14").Value = '25.53'
$ws.Range("E14").Value = '  -2.82%  '
$ws.Range("D15").Value = '61.864.07'
$ws.Range("E15").Value = '  -2.09%  '
$ws.Range("D16").Value = '0.0000144'
$ws.Range("E16").Value = '  -0.13%  '
$ws.Range("D17").Value = '2.583.73'
$ws.Range("E17").Value = '  -3.79%  '
$ws.Range("D18").Value = '11.63'
$ws.Range("E18").Value = '  -3.12%  '
$ws.Range("E19").Value = '  -0.30%  '
$ws.Range("D20").Value = '337.63'
$ws.Range("E20").Value = '  -1.20%  '
$ws.Range("D21").Value = '6.02'
$ws.Range("E21").Value = '  -4.48%  '
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  +0.27%  '
$ws.Range("D23").Value = '0.492'
$ws.Range("E23").Value = '  -2.81%  '
$ws.Range("D24").Value = '63.62'
$ws.Range("E24").Value = '  -0.30%  '
$ws.Range("E25").Value = '  -0.52%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").Value = '8.19'
$ws.Range("E27").Value = '  +0.39%  '
$ws.Range("B28").Value = 'Fetch.AI'
$ws.Range("C28").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D28").Value = '1.36'
$ws.Range("E28").Value = '  +2.94%  '
$ws.Range("B29").Value = 'Aptos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D29").Value = '7.26'
$ws.Range("E29").Value = '  +3.74%  '
$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '0.0₃0836'
$ws.Range("E30").Value = '  -1.77%  '
$ws.Range("E31").Value = '  -2.61%  '
$ws.Range("D32").Value = '162.98'
$ws.Range("E32").Value = '  -1.91%  '
$ws.Range("D33").Value = '4.88'
$ws.Range("E33").Value = '  +2.16%  '
$ws.Range("D35").Value = '19.19'
$ws.Range("E35").Value = '  -1.77%  '
$ws.Range("D36").Value = '1.41'
$ws.Range("E36").Value = '  -1.05%  '
$ws.Range("D37").Value = '1.79'
$ws.Range("E37").Value = '  +1.29%  '
$ws.Range("D38").Value = '328.81'
$ws.Range("E38").Value = '  -2.99%  '
$ws.Range("D40").Value = '0.904'
$ws.Range("E40").Value = '  -4.04%  '
$ws.Range("D41").Value = '3.94'
$ws.Range("E41").Value = '  +0.46%  '
$ws.Range("D42").Value = '37.44'
$ws.Range("E42").Value = '  -1.64%  '
$ws.Range("D43").Value = '20.88'
$ws.Range("E43").Value = '  +0.71%  '
$ws.Range("E44").Value = '  -0.11%  '
$ws.Range("D45").Value = '0.607'
$ws.Range("E45").Value = '  -1.80%  '
$ws.Range("D46").Value = '10.93'
$ws.Range("E46").Value = '  -1.22%  '
$ws.Range("D47").Value = '0.0547'
$ws.Range("E47").Value = '  -2.63%  '
$ws.Range("D48").Value = '2.105.89'
$ws.Range("E48").Value = '  +0.64%  '
$ws.Range("D49").Value = '19.51'
$ws.Range("E49").Value = '  -3.49%  '
$ws.Range("D50").Value = '0.0966'
$ws.Range("E50").Value = '  -0.47%  '
$ws.Range("D51").Value = '0.0238'
$ws.Range("E51").Value = '  -0.87%  '
